$d = $word.ActiveDocument

# Pull the whole package as flat OPC WordprocessingML so we can retarget
# the picture "name" attributes (wp:docPr / pic:cNvPr) that live inside
# the header/footer parts holding the Pearson/BTEC logos.
$xml = $d.WordOpenXML

# --- Footer (first page) logo: id="3" image2.png -> image1.png ---
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image2.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"')

# --- Footer (default) logo: id="2" image2.png -> image1.png ---
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"')

# --- pic:cNvPr copies inside both Pearson logo drawings (id="0") ---
$xml = $xml.Replace(
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"',
    'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"')

# --- Header (first page) BTEC logo: image1.jpg -> image2.jpg ---
$xml = $xml.Replace(
    'descr="BTec_Logo-Orange" id="1" name="image1.jpg"',
    'descr="BTec_Logo-Orange" id="1" name="image2.jpg"')

$xml = $xml.Replace(
    'descr="BTec_Logo-Orange" id="0" name="image1.jpg"',
    'descr="BTec_Logo-Orange" id="0" name="image2.jpg"')

$d.WordOpenXML = $xml
